# Atualização de bases das ligas, do dia: 27-04-2024 às 09:20
# Appends two new match rows (136, 137) to the "India Super League" sheet,
# mirroring the formatting of the last existing data row (135).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last data row onto the two new rows
# so the id column (A) keeps its bold/border style and the date column (D)
# keeps its date number format.
$ws.Range("A135:AB135").Copy()
$ws.Range("A136:AB137").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 136 ---
$ws.Range("A136").Value = 134
$ws.Range("B136").Value = 8120939
$ws.Range("C136").Value = "India Super League"
$ws.Range("D136").Value = 45405.45833333334
$ws.Range("E136").Value = "Odisha FC"
$ws.Range("F136").Value = "Mohun Bagan SG"
$ws.Range("G136").Value = 2
$ws.Range("H136").Value = 1
$ws.Range("I136").Value = "H"
$ws.Range("J136").Value = 3.4
$ws.Range("K136").Value = 3.5
$ws.Range("L136").Value = 2
$ws.Range("M136").Value = 3.1
$ws.Range("N136").Value = 3.4
$ws.Range("O136").Value = 2.15
$ws.Range("P136").Value = 0.25
$ws.Range("Q136").Value = 1.9
$ws.Range("R136").Value = 1.95
$ws.Range("S136").Value = 2.75
$ws.Range("T136").Value = 1.875
$ws.Range("U136").Value = 1.975
$ws.Range("V136").Value = 2.1
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = 0.8999999999999999
$ws.Range("Z136").Value = -1
$ws.Range("AA136").Value = 0.4375
$ws.Range("AB136").Value = -0.5

# --- Row 137 ---
$ws.Range("A137").Value = 135
$ws.Range("B137").Value = 8124823
$ws.Range("C137").Value = "India Super League"
$ws.Range("D137").Value = 45406.45833333334
$ws.Range("E137").Value = "FC Goa"
$ws.Range("F137").Value = "Mumbai City FC"
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = 3
$ws.Range("I137").Value = "A"
$ws.Range("J137").Value = 2.05
$ws.Range("K137").Value = 3.5
$ws.Range("L137").Value = 3.2
$ws.Range("M137").Value = 2.3
$ws.Range("N137").Value = 3.1
$ws.Range("O137").Value = 3
$ws.Range("P137").Value = -0.25
$ws.Range("Q137").Value = 2.025
$ws.Range("R137").Value = 1.825
$ws.Range("S137").Value = 2.5
$ws.Range("T137").Value = 2
$ws.Range("U137").Value = 1.85
$ws.Range("V137").Value = -1
$ws.Range("W137").Value = -1
$ws.Range("X137").Value = 2
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = 0.825
$ws.Range("AA137").Value = 1
$ws.Range("AB137").Value = -1
